$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All target columns (B-E) hold literal text in the source data (coin names,
# URLs, price strings using "." as both decimal and thousands separator, and
# padded percentage strings) -- force text format first so Excel does not
# auto-coerce/round these into numbers (which would e.g. turn "0.990" into 0.99
# or "20.10" into 20.1).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "55.848.74"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.948.33"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.42%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.990"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.94%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "499.67"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.28"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.429"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "7.06"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.106"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +0.60%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.366"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +3.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "3.492.57"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.64%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.02"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000161"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "55.559.10"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -2.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.972.02"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.95"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.11%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.82"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.57"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.43%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.491"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.87"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +1.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.105.14"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.78%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +1.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.160"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0879"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.48"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.18%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.27%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.87%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.10"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.79"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.37%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.56"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.77"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.35%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.25"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.36%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.46"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0657"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.986.22"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.19%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.51"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.23%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.989"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.77"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +2.12%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.156.92"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.67%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.932"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.01%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "Cosmos"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.88"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.21%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "19.66"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.80%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0235"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.58%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0849"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.71%  "
